# Update the Rules worksheet: change the greeting text in E8 and select that cell,
# matching the author's "update file with jgit" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E8").Value = "GIT UPDATE"
$ws.Range("E8").Select()
